# Slide 2 ("计算机网络01" diagram): reposition the "主机" (host) and
# "服务器" (server) rectangles, and add a new "路由器" (router) rectangle
# cloned from the server shape so it keeps the same style/appearance.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- 主机 ("矩形 50") : move to new position -------------------------------
$hostShape = $s.Shapes.Item(23)
$hostShape.Left = 58.729841319685036
$hostShape.Top  = 435.8692169984252

# --- 服务器 ("矩形 51") : move to new position ------------------------------
$server = $s.Shapes.Item(24)
$server.Left = 830.6323548047244
$server.Top  = 440.7726898653543

# --- 路由器 ("矩形 26") : new shape, cloned from the server rectangle so it
#     inherits the same style (fill/line/effect/font refs) and text body
#     formatting, then repositioned and retitled. ---------------------------
$router = $server.Duplicate()
$router.Name = "矩形 26"
$router.Left = 444.68110666220474
$router.Top  = 434.4319000637795
$router.TextFrame.TextRange.Text = "路由器"
